$wb = $excel.ActiveWorkbook

# --- ALC: 14 cell changes ---
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
  @(74, 8, 118424.38),
  @(74, 9, 147109.7),
  @(74, 10, 22806.666),
  @(74, 11, 147109.7),
  @(74, 12, 22806.666),
  @(74, 13, -146173.7),
  @(74, 14, -24678.666),
  @(77, 8, 118424.38),
  @(77, 9, 147109.7),
  @(77, 10, 22806.666),
  @(77, 11, 735548.5),
  @(77, 12, 114033.33),
  @(77, 13, -730868.5),
  @(77, 14, -123393.33)
)
foreach ($e in $edits) {
    if ($null -eq $e[2]) {
        $ws.Cells.Item($e[0], $e[1]).Value = ""
    } else {
        $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
    }
}

# --- BSM: 125 cell changes ---
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
  @(117, 8, 0),
  @(117, 9, 0),
  @(117, 10, 0),
  @(117, 11, 0),
  @(117, 12, 0),
  @(118, 8, 200000),
  @(118, 9, 0),
  @(118, 10, 200000),
  @(118, 11, 0),
  @(118, 12, 200000),
  @(118, 14, -203314),
  @(119, 8, 0),
  @(119, 9, 0),
  @(119, 10, 0),
  @(119, 11, 0),
  @(119, 12, 0),
  @(120, 8, 0),
  @(120, 9, 0),
  @(120, 10, 0),
  @(120, 11, 0),
  @(120, 12, 0),
  @(122, 8, 0),
  @(122, 9, 0),
  @(122, 10, 0),
  @(122, 11, 0),
  @(122, 12, 0),
  @(123, 8, 54995),
  @(123, 9, 0),
  @(123, 10, 54995),
  @(123, 11, 0),
  @(123, 12, 54995),
  @(123, 14, -64795),
  @(124, 8, 0),
  @(124, 9, 0),
  @(124, 10, 0),
  @(124, 11, 0),
  @(124, 12, 0),
  @(125, 8, 0),
  @(125, 9, 0),
  @(125, 10, 0),
  @(125, 11, 0),
  @(125, 12, 0),
  @(126, 8, 0),
  @(126, 9, 0),
  @(126, 10, 0),
  @(126, 11, 0),
  @(126, 12, 0),
  @(127, 8, 64553.4),
  @(127, 9, 44780),
  @(127, 10, 69496.75),
  @(127, 11, 44780),
  @(127, 12, 69496.75),
  @(127, 13, -39820),
  @(127, 14, -79416.75),
  @(128, 8, 0),
  @(128, 9, 0),
  @(128, 10, 0),
  @(128, 11, 0),
  @(128, 12, 0),
  @(129, 8, 0),
  @(129, 9, 0),
  @(129, 10, 0),
  @(129, 11, 0),
  @(129, 12, 0),
  @(130, 8, 0),
  @(130, 9, 0),
  @(130, 10, 0),
  @(130, 11, 0),
  @(130, 12, 0),
  @(131, 8, 99997.664),
  @(131, 9, 0),
  @(131, 10, 99997.664),
  @(131, 11, 0),
  @(131, 12, 99997.664),
  @(131, 14, -110077.664),
  @(132, 8, 74999.5),
  @(132, 9, 70000),
  @(132, 10, 79999),
  @(132, 11, 70000),
  @(132, 12, 79999),
  @(132, 13, -64940),
  @(132, 14, -90119),
  @(133, 8, 0),
  @(133, 9, 0),
  @(133, 10, 0),
  @(133, 11, 0),
  @(133, 12, 0),
  @(134, 8, 303.6),
  @(134, 9, 303.6),
  @(134, 10, 0),
  @(134, 11, 910.8000000000001),
  @(134, 12, 0),
  @(134, 13, 1624.2),
  @(135, 8, 0),
  @(135, 9, 0),
  @(135, 10, 0),
  @(135, 11, 0),
  @(135, 12, 0),
  @(137, 8, 0),
  @(137, 9, 0),
  @(137, 10, 0),
  @(137, 11, 0),
  @(137, 12, 0),
  @(138, 8, 0),
  @(138, 9, 0),
  @(138, 10, 0),
  @(138, 11, 0),
  @(138, 12, 0),
  @(139, 8, 0),
  @(139, 9, 0),
  @(139, 10, 0),
  @(139, 11, 0),
  @(139, 12, 0),
  @(140, 8, 288498.5),
  @(140, 9, 0),
  @(140, 10, 288498.5),
  @(140, 11, 0),
  @(140, 12, 288498.5),
  @(140, 14, -298858.5),
  @(141, 8, 83498.5),
  @(141, 9, 0),
  @(141, 10, 83498.5),
  @(141, 11, 0),
  @(141, 12, 83498.5),
  @(141, 14, -93858.5)
)
foreach ($e in $edits) {
    if ($null -eq $e[2]) {
        $ws.Cells.Item($e[0], $e[1]).Value = ""
    } else {
        $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
    }
}

# --- CRP: 21 cell changes ---
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
  @(99, 8, 873.3333),
  @(99, 9, 869.0909),
  @(99, 10, 880),
  @(99, 11, 869.0909),
  @(99, 12, 880),
  @(99, 13, 628.9091),
  @(99, 14, -3876),
  @(122, 8, 1713.7273),
  @(122, 9, 885.2),
  @(122, 10, 9999),
  @(122, 11, 2655.6),
  @(122, 12, 29997),
  @(122, 13, -205.6000000000004),
  @(122, 14, -34897),
  @(126, 8, 873.3333),
  @(126, 9, 869.0909),
  @(126, 10, 880),
  @(126, 11, 2607.2727),
  @(126, 12, 2640),
  @(126, 13, -137.2727),
  @(126, 14, -7580)
)
foreach ($e in $edits) {
    if ($null -eq $e[2]) {
        $ws.Cells.Item($e[0], $e[1]).Value = ""
    } else {
        $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
    }
}

# --- CUL: 141 cell changes ---
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
  @(23, 8, 107.2),
  @(23, 9, 97.75),
  @(23, 11, 293.25),
  @(23, 13, -58.25),
  @(68, 8, 4436.0625),
  @(68, 10, 4498.467),
  @(68, 12, 13495.401),
  @(68, 14, -15117.401),
  @(71, 8, 4436.0625),
  @(71, 10, 4498.467),
  @(71, 12, 40486.20299999999),
  @(71, 14, -48598.20299999999),
  @(120, 8, $null),
  @(120, 9, $null),
  @(120, 10, $null),
  @(120, 11, $null),
  @(120, 12, $null),
  @(121, 8, $null),
  @(121, 9, $null),
  @(121, 10, $null),
  @(121, 11, $null),
  @(121, 12, $null),
  @(121, 13, $null),
  @(121, 14, $null),
  @(122, 8, $null),
  @(122, 9, $null),
  @(122, 10, $null),
  @(122, 11, $null),
  @(122, 12, $null),
  @(123, 8, $null),
  @(123, 9, $null),
  @(123, 10, $null),
  @(123, 11, $null),
  @(123, 12, $null),
  @(124, 8, $null),
  @(124, 9, $null),
  @(124, 10, $null),
  @(124, 11, $null),
  @(124, 12, $null),
  @(125, 8, $null),
  @(125, 9, $null),
  @(125, 10, $null),
  @(125, 11, $null),
  @(125, 12, $null),
  @(126, 8, $null),
  @(126, 9, $null),
  @(126, 10, $null),
  @(126, 11, $null),
  @(126, 12, $null),
  @(126, 13, $null),
  @(127, 8, $null),
  @(127, 9, $null),
  @(127, 10, $null),
  @(127, 11, $null),
  @(127, 12, $null),
  @(128, 8, $null),
  @(128, 9, $null),
  @(128, 10, $null),
  @(128, 11, $null),
  @(128, 12, $null),
  @(128, 13, $null),
  @(129, 8, $null),
  @(129, 9, $null),
  @(129, 10, $null),
  @(129, 11, $null),
  @(129, 12, $null),
  @(129, 13, $null),
  @(129, 14, $null),
  @(130, 8, $null),
  @(130, 9, $null),
  @(130, 10, $null),
  @(130, 11, $null),
  @(130, 12, $null),
  @(130, 13, $null),
  @(131, 8, $null),
  @(131, 9, $null),
  @(131, 10, $null),
  @(131, 11, $null),
  @(131, 12, $null),
  @(131, 13, $null),
  @(131, 14, $null),
  @(132, 8, $null),
  @(132, 9, $null),
  @(132, 10, $null),
  @(132, 11, $null),
  @(132, 12, $null),
  @(132, 13, $null),
  @(133, 8, $null),
  @(133, 9, $null),
  @(133, 10, $null),
  @(133, 11, $null),
  @(133, 12, $null),
  @(133, 13, $null),
  @(133, 14, $null),
  @(134, 8, $null),
  @(134, 9, $null),
  @(134, 10, $null),
  @(134, 11, $null),
  @(134, 12, $null),
  @(134, 13, $null),
  @(134, 14, $null),
  @(136, 8, $null),
  @(136, 9, $null),
  @(136, 10, $null),
  @(136, 11, $null),
  @(136, 12, $null),
  @(136, 13, $null),
  @(136, 14, $null),
  @(137, 8, $null),
  @(137, 9, $null),
  @(137, 10, $null),
  @(137, 11, $null),
  @(137, 12, $null),
  @(137, 13, $null),
  @(137, 14, $null),
  @(138, 8, $null),
  @(138, 9, $null),
  @(138, 10, $null),
  @(138, 11, $null),
  @(138, 12, $null),
  @(138, 13, $null),
  @(139, 8, $null),
  @(139, 9, $null),
  @(139, 10, $null),
  @(139, 11, $null),
  @(139, 12, $null),
  @(139, 13, $null),
  @(139, 14, $null),
  @(140, 8, $null),
  @(140, 9, $null),
  @(140, 10, $null),
  @(140, 11, $null),
  @(140, 12, $null),
  @(140, 13, $null),
  @(140, 14, $null),
  @(141, 8, $null),
  @(141, 9, $null),
  @(141, 10, $null),
  @(141, 11, $null),
  @(141, 12, $null),
  @(141, 13, $null)
)
foreach ($e in $edits) {
    if ($null -eq $e[2]) {
        $ws.Cells.Item($e[0], $e[1]).Value = ""
    } else {
        $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
    }
}

# --- GSM: 4 cell changes ---
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
  @(49, 8, 0),
  @(49, 10, 0),
  @(49, 12, 0),
  @(49, 14, $null)
)
foreach ($e in $edits) {
    if ($null -eq $e[2]) {
        $ws.Cells.Item($e[0], $e[1]).Value = ""
    } else {
        $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
    }
}

# --- WVR: 7 cell changes ---
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
  @(49, 8, 100022500),
  @(49, 9, 200000000),
  @(49, 10, 44999),
  @(49, 11, 200000000),
  @(49, 12, 44999),
  @(49, 13, -199999770),
  @(49, 14, -45459)
)
foreach ($e in $edits) {
    if ($null -eq $e[2]) {
        $ws.Cells.Item($e[0], $e[1]).Value = ""
    } else {
        $ws.Cells.Item($e[0], $e[1]).Value = $e[2]
    }
}
